$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.715.30"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "'1.629.62"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.81%  "
$ws.Range("D5").Value = "'214.45"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").Value = "'19.46"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.677.71"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.26"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'1.855.10"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "'25.723.56"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "'191.53"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").Value = "'6.25"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("D26").Value = "'141.70"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("D28").Value = "'6.83"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "'15.45"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").Value = "'0.904"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Value = "'1.142.53"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("D39").Value = "'0.540"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").Value = "'100.63"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").Value = "'0.804"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "'1.764.70"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "'0.0508"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "'0.419"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  +4.57%  "
$ws.Range("D51").Value = "'2.34"
$ws.Range("E51").Value = "  -3.45%  "
